$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resource")

# Remove the last data row (row 8: "Mineral" / "자원이 풍부하게 들어있는 광석이다.")
# This reverts the previously-added Mineral resource entry.
$ws.Rows.Item(8).Delete()

# Reset the active selection to A3 (matches the reverted sheet's saved view state)
$ws.Range("A3").Select()
